$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C220").Value = 1047.5175088
$ws.Range("C221").Value = 1058.535933
$ws.Range("C222").Value = 1069.5782869
$ws.Range("C223").Value = 1080.6319802
$ws.Range("C224").Value = 1091.6835496
$ws.Range("C225").Value = 1102.7188188
$ws.Range("C226").Value = 1113.7230871
$ws.Range("C227").Value = 1124.6813447
$ws.Range("C228").Value = 1135.5785107
$ws.Range("C229").Value = 1146.3996876
$ws.Range("C230").Value = 1157.1304256
$ws.Range("C231").Value = 1167.7569866
$ws.Range("C232").Value = 1178.2666007
$ws.Range("C233").Value = 1188.647704
$ws.Range("C234").Value = 1198.890149
$ws.Range("C235").Value = 1208.9853788
$ws.Range("C236").Value = 1218.9265581
$ws.Range("C237").Value = 1228.708656
$ws.Range("C238").Value = 1238.3284772
$ws.Range("C239").Value = 1247.7846407
$ws.Range("C240").Value = 1257.0775085
$ws.Range("C241").Value = 1266.2090669
$ws.Range("C242").Value = 1275.1827675
$ws.Range("C243").Value = 1284.0033345
$ws.Range("C244").Value = 1292.6765456
$ws.Range("C245").Value = 1301.2089974
$ws.Range("C246").Value = 1309.607861
$ws.Range("C247").Value = 1317.8806378
$ws.Range("C248").Value = 1326.0349216
$ws.Range("C249").Value = 1334.0781729
$ws.Range("C250").Value = 1342.0175113
$ws.Range("C251").Value = 1349.8595285
$ws.Range("C252").Value = 1357.6101269
$ws.Range("C253").Value = 1365.2743853
$ws.Range("C254").Value = 1372.8564538
$ws.Range("C255").Value = 1380.3594824
$ws.Range("C256").Value = 1387.785584
$ws.Range("C257").Value = 1395.1358346
$ws.Range("C258").Value = 1402.4103156
$ws.Range("C259").Value = 1409.6081963
$ws.Range("C260").Value = 1416.7278629
$ws.Range("C261").Value = 1423.7670894
$ws.Range("C262").Value = 1430.7232523
$ws.Range("C263").Value = 1437.5935819
$ws.Range("C264").Value = 1444.3754466
$ws.Range("C265").Value = 1451.0666585
$ws.Range("C266").Value = 1457.66579
$ws.Range("C267").Value = 1464.1724887
$ws.Range("C268").Value = 1470.5877758
$ws.Range("C269").Value = 1476.9143126
$ws.Range("C270").Value = 1483.1566233
$ws.Range("C271").Value = 1489.3212593
$ws.Range("C272").Value = 1495.4168946
$ws.Range("C273").Value = 1501.4543441
$ws.Range("C274").Value = 1507.4464959
$ws.Range("C275").Value = 1513.4081532
$ws.Range("C276").Value = 1519.355781
$ws.Range("C277").Value = 1525.3071537
$ws.Range("C278").Value = 1531.2809023
$ws.Range("C279").Value = 1537.2959586
$ws.Range("C280").Value = 1543.3708979
$ws.Range("C281").Value = 1549.5231851
$ws.Range("C282").Value = 1555.7683354
$ws.Range("C283").Value = 1562.1190095
$ws.Range("C284").Value = 1568.5840757
$ws.Range("C285").Value = 1575.1676859
$ws.Range("C286").Value = 1581.8684297
$ws.Range("C287").Value = 1588.678643
$ws.Range("C288").Value = 1595.5839611
$ws.Range("C289").Value = 1602.5632056
$ws.Range("C290").Value = 1609.5886846
$ws.Range("C291").Value = 1616.6269556
$ws.Range("C292").Value = 1623.6400627
$ws.Range("C293").Value = 1630.5871985
$ws.Range("C294").Value = 1637.426687
$ws.Range("C295").Value = 1644.1181261
$ws.Range("C296").Value = 1650.6244966
$ws.Range("C297").Value = 1656.9140361
$ws.Range("C298").Value = 1662.9616992
$ws.Range("C299").Value = 1668.7500769
$ws.Range("C300").Value = 1674.2697188
$ws.Range("C301").Value = 1679.5188728
$ws.Range("C302").Value = 1684.5027252
$ws.Range("C303").Value = 1689.2322657
$ws.Range("C304").Value = 1693.7229275
$ws.Range("C305").Value = 1697.9931473
$ws.Range("C306").Value = 1702.0629745
$ws.Range("C307").Value = 1705.9528234
$ws.Range("C308").Value = 1709.6824334
$ws.Range("C309").Value = 1713.2700639
$ws.Range("C310").Value = 1716.7319282
$ws.Range("C311").Value = 1720.0818497
$ws.Range("C312").Value = 1723.3311083
$ws.Range("C313").Value = 1726.4884454
$ws.Range("C314").Value = 1729.5601887
$ws.Range("C315").Value = 1732.5504637
$ws.Range("C316").Value = 1735.4614647
$ws.Range("C317").Value = 1738.2937585
$ws.Range("C318").Value = 1741.0466056
$ws.Range("C319").Value = 1743.7182826
$ws.Range("C320").Value = 1746.3063985
$ws.Range("C321").Value = 1748.8081947
$ws.Range("C322").Value = 1751.2208265
$ws.Range("C323").Value = 1753.5416192
$ws.Range("C324").Value = 1755.7682959
$ws.Range("C325").Value = 1757.8991721
$ws.Range("C326").Value = 1759.9333145
$ws.Range("C327").Value = 1761.8706581
$ws.Range("C328").Value = 1763.7120799
$ws.Range("C329").Value = 1765.4594256
$ws.Range("C330").Value = 1767.1154903
$ws.Range("C331").Value = 1768.6839515
$ws.Range("C332").Value = 1770.1692605
$ws.Range("C333").Value = 1771.576495
$ws.Range("C334").Value = 1772.9111826
$ws.Range("C335").Value = 1774.1791026
$ws.Range("C336").Value = 1775.3860785
$ws.Range("C337").Value = 1776.5377717
$ws.Range("C338").Value = 1777.6394888
$ws.Range("C339").Value = 1778.6960116
$ws.Range("C340").Value = 1779.7114621
$ws.Range("C341").Value = 1780.6892059
$ws.Range("C342").Value = 1781.6318023
$ws.Range("C343").Value = 1782.5409996
$ws.Range("C344").Value = 1783.4177757
$ws.Range("C345").Value = 1784.2624179
$ws.Range("C347").Value = 1785.8536996
$ws.Range("C349").Value = 1787.3081254
